$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so numeric-looking strings
# (e.g. "63.05", "0.634") are preserved verbatim as text instead of being
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.820.73"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "2.265.72"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "229.87"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").Value = "0.634"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("D7").Value = "63.05"
$ws.Range("E7").Value = "  +3.57%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.446"
$ws.Range("E9").Value = "  +9.63%  "
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  +10.57%  "
$ws.Range("D11").Value = "56.93"
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("D12").Value = "26.11"
$ws.Range("E12").Value = "  +16.68%  "
$ws.Range("E13").Value = "  +1.81%  "
$ws.Range("D14").Value = "2.607.44"
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").Value = "15.57"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("D16").Value = "6.19"
$ws.Range("E16").Value = "  +8.81%  "
$ws.Range("D17").Value = "0.843"
$ws.Range("E17").Value = "  +5.28%  "
$ws.Range("D18").Value = "2.286.02"
$ws.Range("E18").Value = "  +2.13%  "
$ws.Range("D19").Value = "43.762.74"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").Value = "0.0₂01000"
$ws.Range("E20").Value = "  +6.39%  "
$ws.Range("D21").Value = "73.50"
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").Value = "6.06"
$ws.Range("E22").Value = "  -2.15%  "
$ws.Range("D23").Value = "252.67"
$ws.Range("E23").Value = "  +2.84%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("E25").Value = "  -6.95%  "
$ws.Range("E26").Value = "  -2.10%  "
$ws.Range("D27").Value = "3.33"
$ws.Range("E27").Value = "  +24.95%  "
$ws.Range("D28").Value = "10.00"
$ws.Range("E28").Value = "  +2.82%  "
$ws.Range("D29").Value = "171.65"
$ws.Range("E29").Value = "  +1.16%  "
$ws.Range("E30").Value = "  -3.35%  "
$ws.Range("D31").Value = "20.78"
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("D32").Value = "1.38"
$ws.Range("E32").Value = "  -6.14%  "
$ws.Range("D33").Value = "0.124"
$ws.Range("E33").Value = "  +2.92%  "
$ws.Range("E34").Value = "  +6.42%  "
$ws.Range("D35").Value = "4.78"
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("D36").Value = "4.89"
$ws.Range("E36").Value = "  -2.25%  "
$ws.Range("D37").Value = "3.80"
$ws.Range("E37").Value = "  +5.80%  "
$ws.Range("D38").Value = "6.51"
$ws.Range("E38").Value = "  +1.45%  "
$ws.Range("D39").Value = "2.30"
$ws.Range("E39").Value = "  -3.19%  "
$ws.Range("D40").Value = "0.0258"
$ws.Range("E40").Value = "  +3.27%  "
$ws.Range("D41").Value = "0.000240"
$ws.Range("E41").Value = "  +6.92%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "17.64"
$ws.Range("E43").Value = "  +7.63%  "
$ws.Range("D44").Value = "0.0972"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "8.23"
$ws.Range("E45").Value = "  -5.22%  "
$ws.Range("D46").Value = "97.81"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("E47").Value = "  -1.16%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").Value = "4.33"
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("B49").Value = "Celestia"
$ws.Range("C49").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D49").Value = "10.07"
$ws.Range("E49").Value = "  +10.45%  "
$ws.Range("D50").Value = "1.443.06"
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("D51").Value = "2.28"
$ws.Range("E51").Value = "  +1.86%  "

# Restore the original (default) style on column D so no stray number-format
# style is left applied to the cells.
$ws.Range("D2:D51").Style = "Normal"
